$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Helper: read a table cell's text with the trailing cell-mark
# characters (Bell / Cell-Reference + Carriage-Return) stripped.
# ---------------------------------------------------------------
function CellText($tbl, $row, $col) {
    return $tbl.Cell($row, $col).Range.Text.TrimEnd([char]7, [char]13)
}

# ---------------------------------------------------------------
# 1) HMC Information table - simple label re-casing (values unchanged)
# ---------------------------------------------------------------
$d.Content.Find.Execute("HOSTNAME", $false, $false, $false, $false, $false, $true, 1, $false, "Hostname", 2) | Out-Null
$d.Content.Find.Execute("BASE VERSION", $false, $false, $false, $false, $false, $true, 1, $false, "Base Version", 2) | Out-Null
$d.Content.Find.Execute("SERVICE PACK", $false, $false, $false, $false, $false, $true, 1, $false, "Service Pack", 2) | Out-Null
$d.Content.Find.Execute("GATEWAY", $false, $false, $false, $false, $false, $true, 1, $false, "Gateway", 2) | Out-Null
$d.Content.Find.Execute("IP ADDR - ETH0", $false, $false, $false, $false, $false, $true, 1, $false, "Ip Addr - Eth0", 2) | Out-Null
$d.Content.Find.Execute("IP ADDR - ETH1", $false, $false, $false, $false, $false, $true, 1, $false, "Ip Addr - Eth1", 2) | Out-Null
$d.Content.Find.Execute("IP ADDR - ETH2", $false, $false, $false, $false, $false, $true, 1, $false, "Ip Addr - Eth2", 2) | Out-Null
$d.Content.Find.Execute("IP ADDR - ETH3", $false, $false, $false, $false, $false, $true, 1, $false, "Ip Addr - Eth3", 2) | Out-Null

# MODEL / SERIAL appear both in the HMC table and the three Server
# Information tables with identical target casing, so a single
# document-wide replace handles every occurrence.
$d.Content.Find.Execute("MODEL", $false, $false, $false, $false, $false, $true, 1, $false, "Model", 2) | Out-Null
$d.Content.Find.Execute("SERIAL", $false, $false, $false, $false, $false, $true, 1, $false, "Serial", 2) | Out-Null

# ---------------------------------------------------------------
# 2) The nine LPAR Information tables (tables 2-10): collapse the
#    original 11 rows into 6 rows, combining CPU / Virtual
#    Processor / Memory figures into single summary cells.
# ---------------------------------------------------------------
for ($i = 2; $i -le 10; $i++) {
    $tbl = $d.Tables.Item($i)

    $desiredCpu = CellText $tbl 2 2
    $minCpu     = CellText $tbl 3 2
    $maxCpu     = CellText $tbl 4 2
    $desiredVp  = CellText $tbl 5 2
    $minVp      = CellText $tbl 6 2
    $maxVp      = CellText $tbl 7 2
    $entMem     = CellText $tbl 8 2
    $minMem     = CellText $tbl 9 2
    $maxMem     = CellText $tbl 10 2

    $cpuVal = "Desired: " + $desiredCpu + " | Min: " + $minCpu + " | Max: " + $maxCpu
    $vpVal  = "Desired: " + $desiredVp + " | Min: " + $minVp + " | Max: " + $maxVp
    $memVal = "Entitled: " + $entMem + " | Min: " + $minMem + " | Max: " + $maxMem

    $tbl.Cell(1, 1).Range.Text = "Lpar Name"

    $tbl.Cell(2, 1).Range.Text = "CPU"
    $tbl.Cell(2, 2).Range.Text = $cpuVal

    $tbl.Cell(3, 1).Range.Text = "Virtual Processor"
    $tbl.Cell(3, 2).Range.Text = $vpVal

    $tbl.Cell(4, 1).Range.Text = "Memory (Gb)"
    $tbl.Cell(4, 2).Range.Text = $memVal

    # Rows 5-10 (DESIRED/MIN/MAX VIRTUAL PROCESSOR + ENTITLED/MIN/MAX
    # MEMORY) are no longer needed; row 5 is deleted six times so the
    # former POWER SERVER row (originally row 11) slides up to row 5.
    $tbl.Rows.Item(5).Delete()
    $tbl.Rows.Item(5).Delete()
    $tbl.Rows.Item(5).Delete()
    $tbl.Rows.Item(5).Delete()
    $tbl.Rows.Item(5).Delete()
    $tbl.Rows.Item(5).Delete()

    $tbl.Cell(5, 1).Range.Text = "Power Server"
}

# ---------------------------------------------------------------
# 3) System Server Information tables - simple label re-casing
# ---------------------------------------------------------------
$d.Content.Find.Execute("SERVER NAME", $false, $false, $false, $false, $false, $true, 1, $false, "Server Name", 2) | Out-Null
$d.Content.Find.Execute("CPU CORES", $false, $false, $false, $false, $false, $true, 1, $false, "CPU Cores", 2) | Out-Null
$d.Content.Find.Execute("MEMORY (GB)", $false, $false, $false, $false, $false, $true, 1, $false, "Memory (Gb)", 2) | Out-Null
$d.Content.Find.Execute("FIRMWARE LEVEL", $false, $false, $false, $false, $false, $true, 1, $false, "Firmware Level", 2) | Out-Null
$d.Content.Find.Execute("FSP IP ADDRESS", $false, $false, $false, $false, $false, $true, 1, $false, "Fsp Ip Address", 2) | Out-Null
